# Updated symbol list on Wed Jan 18 09:40:34 UTC 2023 with GitHub Actions
# Applies the per-cell text updates to the crypto listing sheet.
# Values are written with a leading apostrophe (quote-prefix) so that
# numeric-looking strings (prices, percentages) are stored as TEXT,
# matching the source workbook's inline-string cells instead of being
# auto-converted to Excel numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '300.93' },
    @{ Cell = 'E2'; Value = '-0.06%' },
    @{ Cell = 'D3'; Value = '32.74' },
    @{ Cell = 'E3'; Value = '3.92%' },
    @{ Cell = 'D4'; Value = '4.950' },
    @{ Cell = 'E4'; Value = '-2.52%' },
    @{ Cell = 'D5'; Value = '0.07704' },
    @{ Cell = 'E5'; Value = '-1.78%' },
    @{ Cell = 'D6'; Value = '1.966' },
    @{ Cell = 'E6'; Value = '-16.07%' },
    @{ Cell = 'D7'; Value = '7.835' },
    @{ Cell = 'E7'; Value = '0.20%' },
    @{ Cell = 'B8'; Value = 'GateToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'D8'; Value = '3.801' },
    @{ Cell = 'E8'; Value = '-0.93%' },
    @{ Cell = 'B9'; Value = 'MXToken' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D9'; Value = '0.9201' },
    @{ Cell = 'E9'; Value = '0.24%' },
    @{ Cell = 'B10'; Value = 'WazirX' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'D10'; Value = '0.1750' },
    @{ Cell = 'E10'; Value = '-0.47%' },
    @{ Cell = 'B11'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D11'; Value = '0.07784' },
    @{ Cell = 'E11'; Value = '2.52%' },
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D12'; Value = '0.08591' },
    @{ Cell = 'E12'; Value = '-7.26%' },
    @{ Cell = 'B13'; Value = 'BitrueCoin' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'D13'; Value = '0.03191' },
    @{ Cell = 'E13'; Value = '6.77%' },
    @{ Cell = 'B14'; Value = 'BitMartToken' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'D14'; Value = '0.1002' },
    @{ Cell = 'E14'; Value = '0.09%' },
    @{ Cell = 'B15'; Value = 'BitForexToken' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'D15'; Value = '0.001518' },
    @{ Cell = 'E15'; Value = '0.53%' },
    @{ Cell = 'B16'; Value = 'TigerCash' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'D16'; Value = '0.005813' },
    @{ Cell = 'E16'; Value = '-0.48%' },
    @{ Cell = 'B17'; Value = 'LEO' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'D17'; Value = '3.460' },
    @{ Cell = 'E17'; Value = '-0.34%' },
    @{ Cell = 'D18'; Value = '2.153' },
    @{ Cell = 'E18'; Value = '-4.21%' },
    @{ Cell = 'E19'; Value = '2.36%' },
    @{ Cell = 'E20'; Value = '-0.12%' },
    @{ Cell = 'D21'; Value = '4.271' },
    @{ Cell = 'E21'; Value = '5.36%' },
    @{ Cell = 'D22'; Value = '0.1993' },
    @{ Cell = 'E22'; Value = '11.36%' },
    @{ Cell = 'D23'; Value = '0.04523' },
    @{ Cell = 'E23'; Value = '-2.01%' },
    @{ Cell = 'D24'; Value = '0.001223' },
    @{ Cell = 'E24'; Value = '-2.18%' },
    @{ Cell = 'D25'; Value = '0.004412' },
    @{ Cell = 'E25'; Value = '-1.40%' },
    @{ Cell = 'D26'; Value = '0.0001252' },
    @{ Cell = 'E26'; Value = '0.17%' },
    @{ Cell = 'D39'; Value = '0.01700' },
    @{ Cell = 'E39'; Value = '-3.68%' },
    @{ Cell = 'D40'; Value = '0.04692' },
    @{ Cell = 'E40'; Value = '-1.86%' },
    @{ Cell = 'D41'; Value = '0.007488' },
    @{ Cell = 'E41'; Value = '3.88%' },
    @{ Cell = 'D42'; Value = '0.1353' },
    @{ Cell = 'E42'; Value = '-0.48%' },
    @{ Cell = 'D43'; Value = '0.002333' },
    @{ Cell = 'E43'; Value = '6.57%' },
    @{ Cell = 'D44'; Value = '0.01053' },
    @{ Cell = 'E44'; Value = '1.93%' },
    @{ Cell = 'D45'; Value = '0.00006261' },
    @{ Cell = 'E45'; Value = '-0.07%' },
    @{ Cell = 'D46'; Value = '0.00000000751' },
    @{ Cell = 'E46'; Value = '0.17%' },
    @{ Cell = 'D47'; Value = '0.8206' },
    @{ Cell = 'E47'; Value = '10.39%' },
    @{ Cell = 'D49'; Value = '0.00002103' },
    @{ Cell = 'E49'; Value = '0.17%' },
    @{ Cell = 'D50'; Value = '0.0002003' },
    @{ Cell = 'E50'; Value = '0.17%' }
)

$quote = [string][char]39

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $quote + $u.Value
}
